$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 15615.7
$ws.Range("I28").Value = 639.3333
$ws.Range("J28").Value = 22034.143
$ws.Range("K28").Value = 639.3333
$ws.Range("L28").Value = 22034.143
$ws.Range("M28").Value = -154.3333
$ws.Range("N28").Value = -23004.143
$ws.Range("H33").Value = 94
$ws.Range("I33").Value = 98.545456
$ws.Range("J33").Value = 69
$ws.Range("K33").Value = 98.545456
$ws.Range("L33").Value = 69
$ws.Range("M33").Value = 130.454544
$ws.Range("N33").Value = -527
$ws.Range("H53").Value = 193.41667
$ws.Range("I53").Value = 310.33334
$ws.Range("J53").Value = 154.44444
$ws.Range("K53").Value = 310.33334
$ws.Range("L53").Value = 154.44444
$ws.Range("M53").Value = 326.66666
$ws.Range("N53").Value = -1428.44444
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()
$ws.Range("H112").Value = 4668.9023
$ws.Range("J112").Value = 5008.5527
$ws.Range("L112").Value = 15025.6581
$ws.Range("N112").Value = -17241.6581
$ws.Range("H137").Value = 1879.9
$ws.Range("I137").Value = 1917.1111
$ws.Range("J137").Value = 1545
$ws.Range("K137").Value = 5751.3333
$ws.Range("L137").Value = 4635
$ws.Range("M137").Value = -3201.3333
$ws.Range("N137").Value = -9735
$ws.Range("H138").Value = 1878.1146
$ws.Range("I138").Value = 1400.1282
$ws.Range("J138").Value = 2205.158
$ws.Range("K138").Value = 4200.3846
$ws.Range("L138").Value = 6615.474
$ws.Range("M138").Value = 939.6153999999997
$ws.Range("N138").Value = -16895.474
$ws.Range("H141").Value = 3026
$ws.Range("I141").Value = 1021.59375
$ws.Range("J141").Value = 7959.923
$ws.Range("K141").Value = 3064.78125
$ws.Range("L141").Value = 23879.769
$ws.Range("M141").Value = 2115.21875
$ws.Range("N141").Value = -34239.769

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 34507
$ws.Range("J7").Value = 37972
$ws.Range("L7").Value = 37972
$ws.Range("N7").Value = -38200
$ws.Range("H32").Value = 801478
$ws.Range("I32").Value = 869733.2
$ws.Range("J32").Value = 27919
$ws.Range("K32").Value = 869733.2
$ws.Range("L32").Value = 27919
$ws.Range("M32").Value = -869446.2
$ws.Range("N32").Value = -28493
$ws.Range("H61").Value = 2880.3428
$ws.Range("I61").Value = 2554.0952
$ws.Range("J61").Value = 3369.7144
$ws.Range("K61").Value = 2554.0952
$ws.Range("L61").Value = 3369.7144
$ws.Range("M61").Value = -2342.0952
$ws.Range("N61").Value = -3793.7144
$ws.Range("H74").Value = 904.95746
$ws.Range("I74").Value = 697.64514
$ws.Range("J74").Value = 1306.625
$ws.Range("K74").Value = 697.64514
$ws.Range("L74").Value = 1306.625
$ws.Range("M74").Value = 176.35486
$ws.Range("N74").Value = -3054.625
$ws.Range("H77").Value = 904.95746
$ws.Range("I77").Value = 697.64514
$ws.Range("J77").Value = 1306.625
$ws.Range("K77").Value = 3488.2257
$ws.Range("L77").Value = 6533.125
$ws.Range("M77").Value = 879.7743
$ws.Range("N77").Value = -15269.125
$ws.Range("H132").Value = 2282.7866
$ws.Range("I132").Value = 1730.4084
$ws.Range("J132").Value = 4461.6113
$ws.Range("K132").Value = 5191.2252
$ws.Range("L132").Value = 13384.8339
$ws.Range("M132").Value = -2661.2252
$ws.Range("N132").Value = -18444.8339
$ws.Range("H136").Value = 2880.3428
$ws.Range("I136").Value = 2554.0952
$ws.Range("J136").Value = 3369.7144
$ws.Range("K136").Value = 7662.285600000001
$ws.Range("L136").Value = 10109.1432
$ws.Range("M136").Value = -5112.285600000001
$ws.Range("N136").Value = -15209.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3523.2144
$ws.Range("I134").Value = 2715.4736
$ws.Range("J134").Value = 5228.4443
$ws.Range("K134").Value = 8146.4208
$ws.Range("L134").Value = 15685.3329
$ws.Range("M134").Value = -5611.4208
$ws.Range("N134").Value = -20755.3329

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3254.519
$ws.Range("I31").Value = 1076.0488
$ws.Range("J31").Value = 5604.9736
$ws.Range("K31").Value = 1076.0488
$ws.Range("L31").Value = 5604.9736
$ws.Range("M31").Value = -781.0488
$ws.Range("N31").Value = -6194.9736
$ws.Range("H34").Value = 3254.519
$ws.Range("I34").Value = 1076.0488
$ws.Range("J34").Value = 5604.9736
$ws.Range("K34").Value = 1076.0488
$ws.Range("L34").Value = 5604.9736
$ws.Range("M34").Value = -874.0488
$ws.Range("N34").Value = -6008.9736
$ws.Range("H58").Value = 1085.0613
$ws.Range("I58").Value = 839.4286
$ws.Range("J58").Value = 1412.5714
$ws.Range("K58").Value = 839.4286
$ws.Range("L58").Value = 1412.5714
$ws.Range("M58").Value = -636.4286
$ws.Range("N58").Value = -1818.5714
$ws.Range("H132").Value = 1438.3829
$ws.Range("I132").Value = 1215.8235
$ws.Range("J132").Value = 2020.4615
$ws.Range("K132").Value = 3647.4705
$ws.Range("L132").Value = 6061.3845
$ws.Range("M132").Value = -1117.4705
$ws.Range("N132").Value = -11121.3845
$ws.Range("H134").Value = 2089.3977
$ws.Range("I134").Value = 2278.5303
$ws.Range("J134").Value = 1355.1177
$ws.Range("K134").Value = 6835.590899999999
$ws.Range("L134").Value = 4065.3531
$ws.Range("M134").Value = -4300.590899999999
$ws.Range("N134").Value = -9135.3531
$ws.Range("H136").Value = 1085.0613
$ws.Range("I136").Value = 839.4286
$ws.Range("J136").Value = 1412.5714
$ws.Range("K136").Value = 2518.2858
$ws.Range("L136").Value = 4237.7142
$ws.Range("M136").Value = 31.71420000000035
$ws.Range("N136").Value = -9337.7142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1901.2273
$ws.Range("J5").Value = 2426.6667
$ws.Range("L5").Value = 7280.000100000001
$ws.Range("N5").Value = -7504.000100000001
$ws.Range("H112").Value = 6675
$ws.Range("I112").Value = 6966.6665
$ws.Range("J112").Value = 6500
$ws.Range("K112").Value = 20899.9995
$ws.Range("L112").Value = 19500
$ws.Range("M112").Value = -19791.9995
$ws.Range("N112").Value = -21716
$ws.Range("H122").Value = 4090.5334
$ws.Range("J122").Value = 9508.416999999999
$ws.Range("L122").Value = 85575.753
$ws.Range("N122").Value = -90475.753
$ws.Range("H132").Value = 2625.6216
$ws.Range("I132").Value = 2391.5
$ws.Range("J132").Value = 2738
$ws.Range("K132").Value = 21523.5
$ws.Range("L132").Value = 24642
$ws.Range("M132").Value = -18993.5
$ws.Range("N132").Value = -29702
$ws.Range("H135").Value = 1901.2273
$ws.Range("J135").Value = 2426.6667
$ws.Range("L135").Value = 21840.0003
$ws.Range("N135").Value = -26910.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 32970.5
$ws.Range("I26").Value = 19798
$ws.Range("J26").Value = 37361.332
$ws.Range("K26").Value = 19798
$ws.Range("L26").Value = 37361.332
$ws.Range("M26").Value = -19518
$ws.Range("N26").Value = -37921.332
$ws.Range("H50").Value = 32970.5
$ws.Range("I50").Value = 19798
$ws.Range("J50").Value = 37361.332
$ws.Range("K50").Value = 19798
$ws.Range("L50").Value = 37361.332
$ws.Range("M50").Value = -19300
$ws.Range("N50").Value = -38357.332
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H113").Value = 79574.234
$ws.Range("I113").Value = 93496.82000000001
$ws.Range("K113").Value = 93496.82000000001
$ws.Range("M113").Value = -91326.82000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1939.5541
$ws.Range("I132").Value = 1680.2
$ws.Range("J132").Value = 2479.875
$ws.Range("K132").Value = 5040.6
$ws.Range("L132").Value = 7439.625
$ws.Range("M132").Value = -2510.6
$ws.Range("N132").Value = -12499.625
$ws.Range("H136").Value = 5557661.5
$ws.Range("I136").Value = 2229.4546
$ws.Range("J136").Value = 20835100
$ws.Range("K136").Value = 6688.3638
$ws.Range("L136").Value = 62505300
$ws.Range("M136").Value = -4138.3638
$ws.Range("N136").Value = -62510400
$ws.Range("H137").Value = 36600
$ws.Range("H138").Value = 41463.77
$ws.Range("I138").Value = 42000
$ws.Range("J138").Value = 41419.082
$ws.Range("K138").Value = 42000
$ws.Range("L138").Value = 41419.082
$ws.Range("M138").Value = -36860
$ws.Range("N138").Value = -51699.082
$ws.Range("H139").Value = 1178721.1
$ws.Range("I139").Value = 9673000
$ws.Range("J139").Value = 46150.6
$ws.Range("K139").Value = 9673000
$ws.Range("L139").Value = 46150.6
$ws.Range("M139").Value = -9667860
$ws.Range("N139").Value = -56430.6
$ws.Range("H140").Value = 41933.375
$ws.Range("J140").Value = 41933.375
$ws.Range("L140").Value = 41933.375
$ws.Range("N140").Value = -52293.375
$ws.Range("H141").Value = 99714
$ws.Range("J141").Value = 99714
$ws.Range("L141").Value = 99714
$ws.Range("N141").Value = -110074

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 759.95654
$ws.Range("I107").Value = 774.8125
$ws.Range("J107").Value = 726
$ws.Range("K107").Value = 2324.4375
$ws.Range("L107").Value = 2178
$ws.Range("M107").Value = -404.4375
$ws.Range("N107").Value = -6018
$ws.Range("H132").Value = 2200.5278
$ws.Range("I132").Value = 2216.2083
$ws.Range("K132").Value = 6648.624899999999
$ws.Range("M132").Value = -4118.624899999999
$ws.Range("H136").Value = 2270.4844
$ws.Range("I136").Value = 1872.8654
$ws.Range("J136").Value = 1355.1177
$ws.Range("K136").Value = 5618.5962
$ws.Range("M136").Value = -3068.5962
